$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 (Bitcoin)
Set-TextValue "D2" "60.931.99"
$ws.Range("E2").Value = "  +0.75%  "

# Row 3 (Ethereum)
Set-TextValue "D3" "2.675.07"
$ws.Range("E3").Value = "  +2.68%  "

# Row 4 (TetherUSD)
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 (BNB)
Set-TextValue "D5" "579.78"
$ws.Range("E5").Value = "  +1.26%  "

# Row 6 (Solana)
Set-TextValue "D6" "145.25"
$ws.Range("E6").Value = "  +1.83%  "

# Row 7 (USDC)
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  +0.26%  "

# Row 8 (XRP)
Set-TextValue "D8" "0.599"
$ws.Range("E8").Value = "  -0.16%  "

# Row 9 (Toncoin)
Set-TextValue "D9" "6.58"
$ws.Range("E9").Value = "  +1.28%  "

# Row 10 (Dogecoin)
$ws.Range("E10").Value = "  +1.62%  "

# Row 11 (Cardano)
Set-TextValue "D11" "0.380"

# Row 12 (TRON)
Set-TextValue "D12" "0.154"
$ws.Range("E12").Value = "  +1.06%  "

# Row 13 (WrappedliquidstakedEther2.0)
Set-TextValue "D13" "3.145.73"
$ws.Range("E13").Value = "  +2.53%  "

# Row 14 (Avalanche)
Set-TextValue "D14" "25.70"
$ws.Range("E14").Value = "  +10.65%  "

# Row 15 (WrappedBTC)
Set-TextValue "D15" "60.925.15"
$ws.Range("E15").Value = "  +0.68%  "

# Row 16 (ShibaInu)
$ws.Range("E16").Value = "  +1.94%  "

# Row 17 (WrappedEther)
Set-TextValue "D17" "2.671.29"
$ws.Range("E17").Value = "  +2.23%  "

# Row 18 (Chainlink)
Set-TextValue "D18" "11.64"
$ws.Range("E18").Value = "  +2.33%  "

# Row 19 (Polkadot)
$ws.Range("E19").Value = "  +1.77%  "

# Row 20 (BitcoinCash)
Set-TextValue "D20" "351.38"
$ws.Range("E20").Value = "  +1.28%  "

# Row 21 (Uniswap)
Set-TextValue "D21" "6.94"
$ws.Range("E21").Value = "  -0.66%  "

# Row 22 (Dai)
Set-TextValue "D22" "0.999"
$ws.Range("E22").Value = "  +0.08%  "

# Row 23 (Polygon)
Set-TextValue "D23" "0.533"
$ws.Range("E23").Value = "  +1.29%  "

# Row 24 (Litecoin)
Set-TextValue "D24" "64.11"
$ws.Range("E24").Value = "  +1.37%  "

# Row 25 (Binance-PegBSC-USD)
Set-TextValue "D25" "0.999"
$ws.Range("E25").Value = "  +0.17%  "

# Row 26 (Kaspa)
$ws.Range("E26").Value = "  +1.85%  "

# Row 27 (InternetComputer(DFINITY))
Set-TextValue "D27" "8.16"
$ws.Range("E27").Value = "  +5.41%  "

# Row 28 (PancakeSwap)
$ws.Range("E28").Value = "  +7.55%  "

# Row 29 (PEPE)
Set-TextValue "D29" "0.0₃0817"
$ws.Range("E29").Value = "  +3.51%  "

# Row 30 (Aptos)
$ws.Range("E30").Value = "  +6.41%  "

# Row 31 (USDe)
Set-TextValue "D31" "0.999"
$ws.Range("E31").Value = "  +0.13%  "

# Row 32 (Monero)
Set-TextValue "D32" "166.72"
$ws.Range("E32").Value = "  +3.16%  "

# Row 33 (EthereumClassic)
$ws.Range("E33").Value = "  +1.96%  "

# Row 34 (Fetch.AI)
$ws.Range("E34").Value = "  +8.28%  "

# Row 35 (NEARProtocol)
$ws.Range("E35").Value = "  +5.69%  "

# Row 36 (ImmutableX)
$ws.Range("E36").Value = "  +8.42%  "

# Row 37 (Stacks)
Set-TextValue "D37" "1.65"
$ws.Range("E37").Value = "  +3.31%  "

# Row 38 (Bittensor)
Set-TextValue "D38" "329.80"
$ws.Range("E38").Value = "  +11.78%  "

# Row 39 (Filecoin)
Set-TextValue "D39" "4.01"
$ws.Range("E39").Value = "  +4.31%  "

# Row 40 (OKB)
Set-TextValue "D40" "38.42"
$ws.Range("E40").Value = "  +1.57%  "

# Row 41 (SuiNetwork)
Set-TextValue "D41" "0.879"
$ws.Range("E41").Value = "  +3.74%  "

# Row 42 (RenderToken)
Set-TextValue "D42" "5.23"
$ws.Range("E42").Value = "  +6.01%  "

# Row 43 (EnergySwap)
Set-TextValue "D43" "20.56"
$ws.Range("E43").Value = "  +4.03%  "

# Row 44 (Aave)
Set-TextValue "D44" "134.12"
$ws.Range("E44").Value = "  -2.37%  "

# Row 45 (Stellar)
$ws.Range("E45").Value = "  +1.89%  "

# Row 46: was Mantle, now Hedera
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D46" "0.0562"
$ws.Range("E46").Value = "  +3.10%  "

# Row 47: was Hedera, now Mantle
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D47" "0.616"
$ws.Range("E47").Value = "  +0.96%  "

# Row 48 (InjectiveProtocol)
Set-TextValue "D48" "20.53"
$ws.Range("E48").Value = "  +3.76%  "

# Row 49 (FirstDigitalUSD)
$ws.Range("E49").Value = "  +0.37%  "

# Row 50 (VeChain)
$ws.Range("E50").Value = "  +3.04%  "

# Row 51 (Maker)
Set-TextValue "D51" "2.130.28"
$ws.Range("E51").Value = "  +5.08%  "
